# Replace the space-separated header names in row 1 with underscore-separated
# names (e.g. "Format Type" -> "Format_Type"). "NARA_Match_Type" in P1 already
# uses underscores, so it is written back unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Format_Type"
$ws.Range("E1").Value = "Format_Standardized_Name"
$ws.Range("F1").Value = "Format_Identification"
$ws.Range("G1").Value = "Format_Name"
$ws.Range("H1").Value = "Format_Version"
$ws.Range("I1").Value = "Registry_Name"
$ws.Range("J1").Value = "Registry_Key"
$ws.Range("K1").Value = "Format_Note"
$ws.Range("L1").Value = "NARA_Format_Name"
$ws.Range("M1").Value = "NARA_PRONOM_URL"
$ws.Range("N1").Value = "NARA_Risk_Level"
$ws.Range("O1").Value = "NARA_Proposed_Preservation_Plan"
$ws.Range("P1").Value = "NARA_Match_Type"

# Move the active selection from L1:P1 to the single cell O2.
$ws.Range("O2").Select()
